$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replace FAUSTO (004893911 / 29672.58) with EDUARDO (004461070 / 30898.42)
# The account-number column holds text that looks numeric (leading zeros), so
# force text formatting before assigning, then clear the format override
# afterwards so no stray style index is left on the cell.
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004461070"
$ws.Cells.Item(5, 1).ClearFormats()
$ws.Cells.Item(5, 2).Value = "EDUARDO"
$ws.Cells.Item(5, 3).Value = 30898.42

# Insert a new row right after the THIAGO row (row 7), pushing GUSTAVO and
# everything below down by one, then populate it with EULER's data.
$ws.Rows("8:8").Insert()
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "004399832"
$ws.Cells.Item(8, 1).ClearFormats()
$ws.Cells.Item(8, 2).Value = "EULER"
$ws.Cells.Item(8, 3).Value = 12193.5
